$d = $word.ActiveDocument

# The existing _GoBack bookmark sits at the end of the "Team Members present"
# paragraph; it moves down to the end of the new "Josh 25%" paragraph, so
# drop it here and re-create it in the right spot as part of the inserted
# fragment below.
$d.Bookmarks("_GoBack").Delete()

# Remove the old third paragraph entirely (the one with the tab + proofErr
# runs and the "hope" grammar-check split); it is being replaced wholesale
# by the rewritten paragraph + new retro paragraphs below.
$old = $d.Paragraphs(3)
$old.Range.Delete()

# Insert the replacement paragraphs (the merged "This sprint..." paragraph, a
# blank paragraph, the "Measurement Criteria:" heading, and the four
# percentage lines, the last one carrying the relocated _GoBack bookmark and
# trailing space) right after the "Team Members present" paragraph.
$d2 = $word.ActiveDocument
$teamPara = $d2.Paragraphs(2)
$insertAt = $d2.Range($teamPara.Range.End, $teamPara.Range.End)

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$fragment = @"
<w:p $w><w:r><w:t xml:space="preserve">This sprint was not as productive as we wanted it to be. The main struggle that we experienced was this sprint had sprint break in the middle. Therefore, we did not get much finished as a group. The main thing we got done was the carousal for the main page. The one thing was hope to focus on in the future sprints are is to stay on task and work through the necessary stories </w:t></w:r></w:p><w:p $w/><w:p $w><w:r><w:t>Measurement Criteria:</w:t></w:r></w:p><w:p $w><w:r><w:t>Katelynn 25%</w:t></w:r></w:p><w:p $w><w:r><w:t>Brian 25%</w:t></w:r></w:p><w:p $w><w:r><w:t>Cheston 25%</w:t></w:r></w:p><w:p $w><w:r><w:t>Josh 25%</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
"@

$insertAt.InsertXML($fragment)
